$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.6852681951507055
$ws1.Range("C2").Value = -0.9444071000789144
$ws1.Range("B3").Value = 0.7040477222343404
$ws1.Range("C3").Value = 0.3703053666688542
$ws1.Range("B4").Value = 0.4206588654705107
$ws1.Range("C4").Value = -0.3542991439712816

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.438513459931377
$ws2.Range("C2").Value = -0.1965484604575288
$ws2.Range("B3").Value = 0.8045302909411134
$ws2.Range("C3").Value = -0.1681934557501857
$ws2.Range("B4").Value = 0.2642770586559482
$ws2.Range("C4").Value = 0.9442755498641071
